$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.3032696666666667
$ws.Range("H2").Value = 0.909809
$ws.Range("I2").Value = 0.150143210583919
$ws.Range("J2").Value = 0.150143210583919
$ws.Range("M2").Value = 3.646930333333334
$ws.Range("N2").Value = 10.940791
$ws.Range("O2").Value = 0.07359174864485112
$ws.Range("P2").Value = 0.07359174864485112
$ws.Range("Q2").Value = 1.106003346546556
$ws.Range("R2").Value = 9.954030118919
$ws.Range("S2").Value = 0.01104930141402271
$ws.Range("T2").Value = 0.01104930141402271

# Row 3
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.3032696666666667
$ws.Range("H3").Value = 0.909809
$ws.Range("I3").Value = 0.150143210583919
$ws.Range("J3").Value = 0.150143210583919
$ws.Range("O3").Value = 0.7702930298336665
$ws.Range("P3").Value = 0.7702930298336667
$ws.Range("Q3").Value = 11.57666021674467
$ws.Range("R3").Value = 104.189941950702
$ws.Range("S3").Value = 0.1156542685896412
$ws.Range("T3").Value = 0.1156542685896412

# Row 4
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.3032696666666667
$ws.Range("H4").Value = 0.909809
$ws.Range("I4").Value = 0.150143210583919
$ws.Range("J4").Value = 0.150143210583919
$ws.Range("M4").Value = 2.694317333333334
$ws.Range("N4").Value = 8.082952000000001
$ws.Range("O4").Value = 0.05436888172824036
$ws.Range("P4").Value = 0.05436888172824036
$ws.Range("Q4").Value = 0.8171047195742223
$ws.Range("R4").Value = 7.353942476168
$ws.Range("S4").Value = 0.008163118458535377
$ws.Range("T4").Value = 0.008163118458535377

# Row 5
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 0.3032696666666667
$ws.Range("H5").Value = 0.909809
$ws.Range("I5").Value = 0.150143210583919
$ws.Range("J5").Value = 0.150143210583919
$ws.Range("M5").Value = 5.042166
$ws.Range("N5").Value = 15.126498
$ws.Range("O5").Value = 0.1017463397932419
$ws.Range("P5").Value = 0.1017463397932419
$ws.Range("Q5").Value = 1.529136002098
$ws.Range("R5").Value = 13.762224018882
$ws.Range("S5").Value = 0.0152765221217197
$ws.Range("T5").Value = 0.0152765221217197

# Row 6
$ws.Range("I6").Value = 0.6441382016790526
$ws.Range("J6").Value = 0.6441382016790526
$ws.Range("M6").Value = 3.646930333333334
$ws.Range("N6").Value = 10.940791
$ws.Range("O6").Value = 0.07359174864485112
$ws.Range("P6").Value = 0.07359174864485112
$ws.Range("Q6").Value = 4.744929883441667
$ws.Range("R6").Value = 42.704368950975
$ws.Range("S6").Value = 0.04740325663051125
$ws.Range("T6").Value = 0.04740325663051125

# Row 7
$ws.Range("I7").Value = 0.6441382016790526
$ws.Range("J7").Value = 0.6441382016790526
$ws.Range("O7").Value = 0.7702930298336665
$ws.Range("P7").Value = 0.7702930298336667
$ws.Range("S7").Value = 0.4961751670029668
$ws.Range("T7").Value = 0.4961751670029668

# Row 8
$ws.Range("I8").Value = 0.6441382016790526
$ws.Range("J8").Value = 0.6441382016790526
$ws.Range("M8").Value = 2.694317333333334
$ws.Range("N8").Value = 8.082952000000001
$ws.Range("O8").Value = 0.05436888172824036
$ws.Range("P8").Value = 0.05436888172824036
$ws.Range("Q8").Value = 3.505508924466667
$ws.Range("R8").Value = 31.5495803202
$ws.Range("S8").Value = 0.03502107370372985
$ws.Range("T8").Value = 0.03502107370372985

# Row 9
$ws.Range("I9").Value = 0.6441382016790526
$ws.Range("J9").Value = 0.6441382016790526
$ws.Range("M9").Value = 5.042166
$ws.Range("N9").Value = 15.126498
$ws.Range("O9").Value = 0.1017463397932419
$ws.Range("P9").Value = 0.1017463397932419
$ws.Range("Q9").Value = 6.56023612845
$ws.Range("R9").Value = 59.04212515605
$ws.Range("S9").Value = 0.06553870434184468
$ws.Range("T9").Value = 0.06553870434184468

# Row 10
$ws.Range("G10").Value = 0.4155246666666666
$ws.Range("H10").Value = 1.246574
$ws.Range("I10").Value = 0.2057185877370285
$ws.Range("J10").Value = 0.2057185877370285
$ws.Range("M10").Value = 3.646930333333334
$ws.Range("N10").Value = 10.940791
$ws.Range("O10").Value = 0.07359174864485112
$ws.Range("P10").Value = 0.07359174864485112
$ws.Range("Q10").Value = 1.515389511114889
$ws.Range("R10").Value = 13.638505600034
$ws.Range("S10").Value = 0.01513919060031715
$ws.Range("T10").Value = 0.01513919060031715

# Row 11
$ws.Range("G11").Value = 0.4155246666666666
$ws.Range("H11").Value = 1.246574
$ws.Range("I11").Value = 0.2057185877370285
$ws.Range("J11").Value = 0.2057185877370285
$ws.Range("O11").Value = 0.7702930298336665
$ws.Range("P11").Value = 0.7702930298336667
$ws.Range("Q11").Value = 15.86175079937466
$ws.Range("R11").Value = 142.755757194372
$ws.Range("S11").Value = 0.1584635942410586
$ws.Range("T11").Value = 0.1584635942410587

# Row 12
$ws.Range("G12").Value = 0.4155246666666666
$ws.Range("H12").Value = 1.246574
$ws.Range("I12").Value = 0.2057185877370285
$ws.Range("J12").Value = 0.2057185877370285
$ws.Range("M12").Value = 2.694317333333334
$ws.Range("N12").Value = 8.082952000000001
$ws.Range("O12").Value = 0.05436888172824036
$ws.Range("P12").Value = 0.05436888172824036
$ws.Range("Q12").Value = 1.119555311827555
$ws.Range("R12").Value = 10.075997806448
$ws.Range("S12").Value = 0.01118468956597514
$ws.Range("T12").Value = 0.01118468956597514

# Row 13
$ws.Range("G13").Value = 0.4155246666666666
$ws.Range("H13").Value = 1.246574
$ws.Range("I13").Value = 0.2057185877370285
$ws.Range("J13").Value = 0.2057185877370285
$ws.Range("M13").Value = 5.042166
$ws.Range("N13").Value = 15.126498
$ws.Range("O13").Value = 0.1017463397932419
$ws.Range("P13").Value = 0.1017463397932419
$ws.Range("Q13").Value = 2.095144346428
$ws.Range("R13").Value = 18.856299117852
$ws.Range("S13").Value = 0.02093111332967755
$ws.Range("T13").Value = 0.02093111332967755
